$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CasesTab row (B2): the Neo4j/Cypher query text was edited to drop the
# trailing `Cohort` column (co.cohort_description) that the query no longer
# returns ("updated ubc2 10 scripts" per the commit message). Re-typing the
# cell's value causes the old (now-unused) shared string to be dropped and
# the new text to be appended to the shared-string table, which is exactly
# the reordering the sheet data reflects (SamplesTab/FilesTab rows shift
# down one shared-string slot).
$newCasesTabQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC02'] and diag.stage_of_disease in [ 'T2N0M0', 'T2N0M1', 'T2N1M0', 'T2N1M1', 'T3N0M0', 'T3N1M0', 'T3N1M1', 'Not Applicable'] OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@

$ws.Range("B2").Value = $newCasesTabQuery

# --- Row heights: with the shorter query text the wrapped rows re-measure
# to a shorter, common height in the newer Excel build that saved this file.
$ws.Rows.Item(2).RowHeight = 304.5
$ws.Rows.Item(3).RowHeight = 304.5
$ws.Rows.Item(4).RowHeight = 304.5

# --- Selection: the author ended up with B2 selected/in view after the edit.
[void]$ws.Range("B2").Select()
